$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6250373.5
$ws.Range("J17").Value = 6250373.5
$ws.Range("L17").Value = 18751120.5
$ws.Range("N17").Value = -18751456.5

$ws.Range("H19").Value = 895.0714
$ws.Range("I19").Value = 678
$ws.Range("J19").Value = 1015.6667
$ws.Range("K19").Value = 678
$ws.Range("L19").Value = 1015.6667
$ws.Range("M19").Value = -503
$ws.Range("N19").Value = -1365.6667

$ws.Range("H106").Value = 2758.3914
$ws.Range("I106").Value = 2758.3914
$ws.Range("K106").Value = 2758.3914
$ws.Range("M106").Value = -2127.3914

$ws.Range("H135").Value = 1186.8572
$ws.Range("J135").Value = 1800
$ws.Range("L135").Value = 16200
$ws.Range("N135").Value = -21270

$ws.Range("H137").Value = 1898.8334
$ws.Range("I137").Value = 1776.7778
$ws.Range("J137").Value = 2265
$ws.Range("K137").Value = 5330.3334
$ws.Range("L137").Value = 6795
$ws.Range("M137").Value = -2780.3334
$ws.Range("N137").Value = -11895

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2583.7368
$ws.Range("I2").Value = 2281.6365
$ws.Range("J2").Value = 2999.125
$ws.Range("K2").Value = 2281.6365
$ws.Range("L2").Value = 2999.125
$ws.Range("M2").Value = -2168.6365
$ws.Range("N2").Value = -3225.125

$ws.Range("H61").Value = 2042.0217
$ws.Range("I61").Value = 1341.1428
$ws.Range("K61").Value = 1341.1428
$ws.Range("M61").Value = -1129.1428

$ws.Range("H74").Value = 1365.3158
$ws.Range("I74").Value = 1365.3158
$ws.Range("K74").Value = 1365.3158
$ws.Range("M74").Value = -491.3158000000001

$ws.Range("H77").Value = 1365.3158
$ws.Range("I77").Value = 1365.3158
$ws.Range("K77").Value = 6826.579000000001
$ws.Range("M77").Value = -2458.579000000001

$ws.Range("H116").Value = 2583.7368
$ws.Range("I116").Value = 2281.6365
$ws.Range("J116").Value = 2999.125
$ws.Range("K116").Value = 2281.6365
$ws.Range("L116").Value = 2999.125
$ws.Range("M116").Value = 12.36349999999993
$ws.Range("N116").Value = -7587.125

$ws.Range("H136").Value = 2042.0217
$ws.Range("I136").Value = 1341.1428
$ws.Range("K136").Value = 4023.4284
$ws.Range("M136").Value = -1473.4284

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2583.7368
$ws.Range("I3").Value = 2281.6365
$ws.Range("J3").Value = 2999.125
$ws.Range("K3").Value = 2281.6365
$ws.Range("L3").Value = 2999.125
$ws.Range("M3").Value = -2167.6365
$ws.Range("N3").Value = -3227.125

$ws.Range("H99").Value = 3160.8572
$ws.Range("I99").Value = 2917.3333
$ws.Range("K99").Value = 2917.3333
$ws.Range("M99").Value = -1419.3333

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 18573.5
$ws.Range("I62").Value = 19999.572
$ws.Range("J62").Value = 16577
$ws.Range("K62").Value = 19999.572
$ws.Range("L62").Value = 16577
$ws.Range("M62").Value = -19375.572
$ws.Range("N62").Value = -17825

$ws.Range("H65").Value = 18573.5
$ws.Range("I65").Value = 19999.572
$ws.Range("J65").Value = 16577
$ws.Range("K65").Value = 99997.86
$ws.Range("L65").Value = 82885
$ws.Range("M65").Value = -96877.86
$ws.Range("N65").Value = -89125

$ws.Range("H86").Value = 7865.5
$ws.Range("I86").Value = 5154.3335
$ws.Range("K86").Value = 5154.3335
$ws.Range("M86").Value = -4031.3335

$ws.Range("H89").Value = 7865.5
$ws.Range("I89").Value = 5154.3335
$ws.Range("K89").Value = 25771.6675
$ws.Range("M89").Value = -20155.6675

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 870.4286
$ws.Range("I12").Value = 1139.2
$ws.Range("J12").Value = 198.5
$ws.Range("K12").Value = 3417.6
$ws.Range("L12").Value = 595.5
$ws.Range("M12").Value = -3244.6
$ws.Range("N12").Value = -941.5

$ws.Range("H18").Value = 975
$ws.Range("I18").Value = 975
$ws.Range("K18").Value = 2925
$ws.Range("M18").Value = -2756

$ws.Range("H80").Value = 38000
$ws.Range("I80").Value = 38000
$ws.Range("K80").Value = 114000
$ws.Range("M80").Value = -113064

$ws.Range("H83").Value = 38000
$ws.Range("I83").Value = 38000
$ws.Range("K83").Value = 342000
$ws.Range("M83").Value = -337320

$ws.Range("H139").Value = 4902.353
$ws.Range("I139").Value = 3796.3333
$ws.Range("J139").Value = 7556.8
$ws.Range("K139").Value = 11388.9999
$ws.Range("L139").Value = 22670.4
$ws.Range("M139").Value = -6248.999899999999
$ws.Range("N139").Value = -32950.4

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10941.857
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10941.857
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10941.857
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -11481.857

$ws.Range("H73").Value = 10941.857
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10941.857
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10941.857
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -12813.857

$ws.Range("H80").Value = 5661.737
$ws.Range("J80").Value = 10619.5
$ws.Range("L80").Value = 10619.5
$ws.Range("N80").Value = -12615.5

$ws.Range("H83").Value = 5661.737
$ws.Range("J83").Value = 10619.5
$ws.Range("L83").Value = 53097.5
$ws.Range("N83").Value = -63081.5

$ws.Range("H122").Value = 2200
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150

$ws.Range("H123").Value = 45859.57
$ws.Range("I123").Value = 40295.668
$ws.Range("J123").Value = 50032.5
$ws.Range("K123").Value = 40295.668
$ws.Range("L123").Value = 50032.5
$ws.Range("M123").Value = -37845.668
$ws.Range("N123").Value = -54932.5

$ws.Range("H126").Value = 3721.08
$ws.Range("I126").Value = 2309.8462
$ws.Range("J126").Value = 5249.9165
$ws.Range("K126").Value = 6929.5386
$ws.Range("L126").Value = 15749.7495
$ws.Range("M126").Value = -4459.5386
$ws.Range("N126").Value = -20689.7495

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3725.5715
$ws.Range("I68").Value = 2741.818
$ws.Range("K68").Value = 2741.818
$ws.Range("M68").Value = -1992.818

$ws.Range("H71").Value = 3725.5715
$ws.Range("I71").Value = 2741.818
$ws.Range("K71").Value = 13709.09
$ws.Range("M71").Value = -9965.09

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H100").Value = 9464.4
$ws.Range("I100").Value = 9164.166999999999
$ws.Range("K100").Value = 9164.166999999999
$ws.Range("M100").Value = -8623.166999999999

$ws.Range("H122").Value = 11000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 33000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -37900

$ws.Range("H132").Value = 5119.2856
$ws.Range("I132").Value = 3185.5386
$ws.Range("K132").Value = 9556.6158
$ws.Range("M132").Value = -7026.6158

$ws.Range("H136").Value = 4467.0415
$ws.Range("I136").Value = 2634.25
$ws.Range("K136").Value = 7902.75
$ws.Range("M136").Value = -5352.75

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5935.2
$ws.Range("I122").Value = 11185
$ws.Range("J122").Value = 2435.3333
$ws.Range("K122").Value = 33555
$ws.Range("L122").Value = 7305.999899999999
$ws.Range("M122").Value = -31105
$ws.Range("N122").Value = -12205.9999
